# Update the cached text of the datetimeFigureOut date field
# (Date Placeholder) from "ו'/טבת/תשפ"א" to "כ"ב/טבת/תשפ"א" across
# the slide master and every slide layout.

$p = $ppt.ActivePresentation
$newDate = 'כ"ב/טבת/תשפ"א'

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Slide layouts
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
